$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G is labeled "K" (strikeouts -> K), update per regenerated save_data
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 1
$ws.Range("G14").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("G16").Value = 1
